$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values must be swapped between row 23 and row 24:
# A (Id), B (Taxonsorteringsordning), E (TaxonId), F (Artnamn),
# G (Vetenskapligt namn), H (Auktor), Q (Ost), R (Nord)
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $cell23 = $ws.Range($col + "23")
    $cell24 = $ws.Range($col + "24")
    $tmp = $cell23.Value2
    $cell23.Value = $cell24.Value2
    $cell24.Value = $tmp
}
